# Re-sort the detected-rows table (A5:B19) by the (newly regenerated) UUID
# file name in column A, keeping each person's name paired with their row.
# This mirrors a re-run of the face/row detection heuristics (font-size,
# color based) that produced new random UUID prefixes for the _fila_N.png
# files, which in turn changes the alphabetical sort order of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("02377c58-c63b-4909-ae81-26307cbed822_fila_3.png",  "Bryan Alexander Cano"),
    @("04c2c643-8418-4c59-9388-c1fe7e051078_fila_15.png", "Orlando Mauricio Guevara"),
    @("301feecf-684e-4a1f-8cfc-e9df5b86810d_fila_6.png",  "Ronier Jose Rivera"),
    @("426d1ee3-4fdd-400d-91b1-bf0f1760517a_fila_11.png", "Eduardo Domingo Zeledon Merca"),
    @("5adf0760-4db8-4070-a8db-fa3511386cd9_fila_13.png", "Marlon Josue Gonzales Cano"),
    @("861de668-39d4-4f42-848a-3a971aede11c_fila_7.png",  "Cristina Jozabed Carvajal"),
    @("8761e7cc-4053-4584-8484-a93f3eb8f579_fila_1.png",  "Hotep Antonio Ruiz Lezama"),
    @("92d91427-2627-495a-82fc-41c4b429288d_fila_9.png",  "Roman Alfonso Grios Boza"),
    @("ac476561-3453-4971-a125-3dc3dd45b57e_fila_8.png",  "David Orlando Mena Valverd"),
    @("ba29fedb-0bbb-4116-b897-21633ceec33f_fila_12.png", "José Danilo Suárez"),
    @("bf716bf7-4eeb-4e59-9ce7-a0a2cf5d9e07_fila_4.png",  "Yadder Fernando Torres"),
    @("c3373abb-5514-4b69-a93c-d0d5e6b7721c_fila_14.png", "Angel Isaac Alvarez Quiñonez"),
    @("d0821fe8-c624-4bde-bc4a-31285f86fe5f_fila_5.png",  "Erick Espinoza"),
    @("d67424cf-7e61-4dd8-b055-56013819a37a_fila_10.png", "Abraham Silva Ampre"),
    @("f483cc6d-c37b-4d63-b7ed-204354164760_fila_2.png",  "Isabella Dompe Estrada")
)

$startRow = 5
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
